# "Now manual setting for importing support entering column names!"
#
# The "Temperatura nastawienia" column (column I) is no longer populated by
# the manual-entry sheet - that data now comes from the (new) manual column
# naming feature instead, so the header and all of its values are cleared
# out here (formatting/styles stay untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the "Temperatura nastawienia" header (I1) and all of the values
# below it (I2:I29), keeping cell styles intact.
$ws.Range("I1:I29").ClearContents()

# Reflect the last active cell/selection recorded for the sheet.
$ws.Range("J8").Select()
